# Ran code for averaged intensities on spiral schemes.
# This re-runs the averaging computation for ALL orientation schemes,
# including the three new spiral sampling schemes, and writes the
# resulting 18-scheme table (was 15 schemes) back into the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row labels (column B) for data rows 10-19, in their new order.
# Rows 2-9 (HKL header, ND/RD/TD Single, Morris, 3x Ring Perpendicular) are unchanged.
$labels = @{
    10 = "Gaussian-Quadrature"
    11 = "Spiral-90deg-10rot-5space"
    12 = "Spiral-90deg-15rot-5space"
    13 = "Spiral-90deg-10rot-3space"
    14 = "NoRotation-tilt60deg"
    15 = "Rotation-NoTilt"
    16 = "Rotation-60detTilt"
    17 = "HexGrid-90degTilt5degRes"
    18 = "HexGrid-90degTilt22p5degRes"
    19 = "HexGrid-60degTilt5degRes"
}

# New A-column (index) values for data rows 10-19
$indices = @{
    10 = 8
    11 = 9
    12 = 10
    13 = 11
    14 = 12
    15 = 13
    16 = 14
    17 = 15
    18 = 16
    19 = 17
}

# New C:M averaged-intensity data for rows 10-19 (11 values each, columns C..M)
$data = @{
    10 = @(1.732273679167148, 1.686442167232839, 0.7414874669677609, 1.732273679167148, 0.9880047329540944, 1.071215803390088, 0.8604059345742562, 1.686442167232839, 1.2139648171003, 1.473119248133724, 1.179971630714365)
    11 = @(0.8775160230008457, 1.647741017954402, 1.140059471635869, 0.8775160230008457, 0.5938586813806822, 2.583175435397902, 0.868989607860985, 1.647741017954402, 1.393900244795136, 1.135708133897991, 1.285223372871781)
    12 = @(0.8756978369841125, 1.651696757827705, 1.139644025183083, 0.8756978369841125, 0.5955545344539952, 2.578670456575261, 0.8681268382183878, 1.651696757827705, 1.395670391505394, 1.135684114244753, 1.284898408207091)
    13 = @(0.8768895849716551, 1.648122290621395, 1.140056396954657, 0.8768895849716551, 0.5943733154820969, 2.581797812965256, 0.868462399832437, 1.648122290621395, 1.394089343788026, 1.135489464379841, 1.284950300137916)
    14 = @(0.3046320000000002, 2.737640000000001, 1.306147999999999, 0.3046320000000002, 0.8309159999999988, 2.474299999999999, 0.8105199999999971, 2.737640000000001, 2.021894, 1.163263, 1.410692666666666)
    15 = @(0, 3.027662500000005, 1.473899999999997, 0, 0.843062500000001, 2.765837500000004, 0.8463874999999991, 3.027662500000005, 2.250781250000001, 1.125390625000001, 1.492808333333334)
    16 = @(0.4133131813888014, 2.140809199616004, 1.257639591935999, 0.4133131813888014, 0.9186980947967979, 1.981674290073604, 0.9048020599808049, 2.140809199616004, 1.699224395776001, 1.056268788582401, 1.269489402965335)
    17 = @(0.9818658657972137, 0.9973954634041032, 0.9968082978711924, 0.9818658657972137, 0.9918081981030068, 0.9947990907591019, 0.9929500757332735, 0.9973954634041032, 0.9971018806376478, 0.9894838732174307, 0.9926044986113153)
    18 = @(1.190184492971753, 0.9466316812117093, 0.9965097633198127, 1.190184492971753, 0.9683137724766275, 1.030006979687314, 0.9375822955243313, 0.9466316812117093, 0.971570722265761, 1.080877607618757, 1.011538164198591)
    19 = @(1.023928482688378, 0.9518127726826778, 0.9877579946533892, 1.023928482688378, 0.9703618010410472, 0.9818261672495591, 0.997369526484688, 0.9518127726826778, 0.9697853836680335, 0.9968569331782057, 0.9855094574666232)
}

$cols = @("C", "D", "E", "F", "G", "H", "I", "J", "K", "L", "M")

for ($row = 10; $row -le 19; $row++) {
    $ws.Cells.Item($row, 1).Value = $indices[$row]
    $ws.Cells.Item($row, 2).Value = $labels[$row]

    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($row, $i + 3).Value = $vals[$i]
    }
}

# Rows 17-19 are brand new rows; give column A the same bold/centered/bordered
# format used by the existing index column (A2:A16, and header row B1:M1) by
# copying the formatting from an already-styled cell (A16).
$fmtSource = $ws.Cells.Item(16, 1)
$fmtSource.Copy()
$fmtTarget = $ws.Range($ws.Cells.Item(17, 1), $ws.Cells.Item(19, 1))
$fmtTarget.PasteSpecial(-4122)
